$wb = $excel.ActiveWorkbook

# --- About sheet ---
$about = $wb.Worksheets.Item("About")
$about.Activate()
# "Last updated" date moved from 2024-01-03 to 2024-03-28
$about.Range("C1").Value = 45379
# View scrolled so row 6 is at the top, without moving the active cell
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

# --- FPIEBP sheet ---
$fpiebp = $wb.Worksheets.Item("FPIEBP")
$fpiebp.Activate()
# Hard coal balancing priorities reordered: production 1 (was 3), imports 3 (was 2), exports 2 (was 1)
$fpiebp.Range("B3").Value = 1
$fpiebp.Range("C3").Value = 3
$fpiebp.Range("D3").Value = 2
# Active cell/selection moved to E3
[void]$fpiebp.Range("E3").Select()
